# Daily attendance processing - rotate the "Recorded By" (column G) list of
# recorders for every data row: move the first name in the comma-separated
# list to the end of the list, unless the list already starts with "System".
#
# Example: "backup@backdoor.com, System, system" -> "System, system, backup@backdoor.com"
#          "dnasr281@gmail.com, System"           -> "System, dnasr281@gmail.com"
#          "System, admin@admin.com"              -> unchanged (already starts with "System")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $text = $cell.Value2

    if ($null -eq $text -or $text -eq "") {
        continue
    }

    $value = "$text"

    if ($value.StartsWith("System")) {
        continue
    }

    if ($value.Contains(",")) {
        $parts = $value -split ",\s*"
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $newValue = [string]::Join(", ", $rotated)
            $cell.Value = $newValue
        }
    }
}
